# Auto-generated PowerShell Excel COM-interop script
# changed MP time limit and corrected error in fixed recourse data
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: per-instance objective/solve time corrections (col B, C) ---
$ws1.Range("B2").Value = -411.5189726818494
$ws1.Range("C2").Value = 41.035422216
$ws1.Range("B3").Value = -411.39604929747
$ws1.Range("C3").Value = 37.704386631
$ws1.Range("B4").Value = -412.34762188775494
$ws1.Range("C4").Value = 32.774448534
$ws1.Range("B5").Value = -416.71129461891667
$ws1.Range("C5").Value = 29.938049389
$ws1.Range("B6").Value = -408.2905898186462
$ws1.Range("C6").Value = 20.644993133
$ws1.Range("B7").Value = -403.9902935908759
$ws1.Range("C7").Value = 27.442563988
$ws1.Range("B8").Value = -400.1919964156508
$ws1.Range("C8").Value = 28.369221222
$ws1.Range("B9").Value = -412.09265747622965
$ws1.Range("C9").Value = 27.04728668
$ws1.Range("B10").Value = -408.4459826348807
$ws1.Range("C10").Value = 28.389731285
$ws1.Range("B11").Value = -403.1039288610207
$ws1.Range("C11").Value = 28.719738159

# --- Sheet1: model-size columns now reflect 50x scaled (multicut) instances ---
$ws1.Range("F2:F11").Value = 50
$ws1.Range("G2:G11").Value = 27700
$ws1.Range("H2:H11").Value = 30250
$ws1.Range("I2:I11").Value = 2500

# --- Per-iteration detail sheets (tabs '1'..'10') ---
$ws = $wb.Worksheets.Item("1")
$ws.Range("D2").Value = 0.9010757257277832
$ws.Range("E2").Value = 84.92934
$ws.Range("B3").Value = -411.5189726818494
$ws.Range("C3").Value = 0.0837459263334878
$ws.Range("D3").Value = 29.604633740277098

$ws = $wb.Worksheets.Item("2")
$ws.Range("D2").Value = 0.03456503136352539
$ws.Range("E2").Value = 81.97293
$ws.Range("B3").Value = -411.39604929747
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 36.7756628814342

$ws = $wb.Worksheets.Item("3")
$ws.Range("D2").Value = 0.07678334272521972
$ws.Range("E2").Value = 86.44331
$ws.Range("B3").Value = -412.34762188775494
$ws.Range("C3").Value = 0.09407025148614129
$ws.Range("D3").Value = 31.316211326493775

$ws = $wb.Worksheets.Item("4")
$ws.Range("D2").Value = 0.0340887926751709
$ws.Range("E2").Value = 85.55874
$ws.Range("B3").Value = -416.71129461891667
$ws.Range("C3").Value = 0.09880801886621136
$ws.Range("D3").Value = 29.069990408236816

$ws = $wb.Worksheets.Item("5")
$ws.Range("D2").Value = 0.11503299599829102
$ws.Range("E2").Value = 84.60907
$ws.Range("B3").Value = -408.2905898186462
$ws.Range("D3").Value = 19.580634745200562

$ws = $wb.Worksheets.Item("6")
$ws.Range("D2").Value = 0.05615988780419922
$ws.Range("E2").Value = 89.81808
$ws.Range("B3").Value = -403.9902935908759
$ws.Range("D3").Value = 26.55329123092859

$ws = $wb.Worksheets.Item("7")
$ws.Range("D2").Value = 0.06278624963916016
$ws.Range("E2").Value = 83.52166
$ws.Range("B3").Value = -400.1919964156508
$ws.Range("C3").Value = 0.00005271016887334036
$ws.Range("D3").Value = 27.271777710144896

$ws = $wb.Worksheets.Item("8")
$ws.Range("D2").Value = 0.05939362430908203
$ws.Range("E2").Value = 85.27605
$ws.Range("B3").Value = -412.09265747622965
$ws.Range("C3").Value = 0.09772822620371914
$ws.Range("D3").Value = 26.089015479168822

$ws = $wb.Worksheets.Item("9")
$ws.Range("D2").Value = 0.028937167536499022
$ws.Range("E2").Value = 82.46345
$ws.Range("B3").Value = -408.4459826348807
$ws.Range("C3").Value = 0.04149298119528847
$ws.Range("D3").Value = 27.672336497333617

$ws = $wb.Worksheets.Item("10")
$ws.Range("D2").Value = 0.03128311739868164
$ws.Range("E2").Value = 84.38256
$ws.Range("B3").Value = -403.1039288610207
$ws.Range("C3").Value = 0.09010305099909127
$ws.Range("D3").Value = 27.802565506419434

Write-Host "Edit applied successfully"
